$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in cell A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 07:52"

# Insert the updated "Hungria" entry into its new sorted position (by Casos totales desc),
# which pushes the old Kuwait row down to row 63 and the old Kazajistan row down to row 64.
$ws.Cells.Item(62, 1).Value = "Hungria"
$ws.Cells.Item(62, 2).Value = 2284
$ws.Cells.Item(62, 3).Value = 116
$ws.Cells.Item(62, 4).Value = 390
$ws.Cells.Item(62, 5).Value = 1655
$ws.Cells.Item(62, 6).Value = 61
$ws.Cells.Item(62, 7).Value = 14
$ws.Cells.Item(62, 8).Value = 239

$ws.Cells.Item(63, 1).Value = "Kuwait"
$ws.Cells.Item(63, 2).Value = 2248
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(63, 4).Value = 443
$ws.Cells.Item(63, 5).Value = 1792
$ws.Cells.Item(63, 6).Value = 50
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 13

$ws.Cells.Item(64, 1).Value = "Kazajistan"
$ws.Cells.Item(64, 2).Value = 2191
$ws.Cells.Item(64, 3).Value = 56
$ws.Cells.Item(64, 4).Value = 515
$ws.Cells.Item(64, 5).Value = 1657
$ws.Cells.Item(64, 6).Value = 29
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 19
